$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73; this shifts the existing rows 73-125
# down to 74-126, automatically carrying all of their data/formatting
# along (matching the diff, which is effectively a "new weekly record"
# prepended to this block and the old last record duplicated at the end).
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new weekly record.
$ws.Cells.Item(73, 1).Value = 7
$ws.Cells.Item(73, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(73, 3).Value = "Ñuble"
$ws.Cells.Item(73, 4).Value = 44589
$ws.Cells.Item(73, 5).Value = 16
$ws.Cells.Item(73, 6).Value = 100112045
$ws.Cells.Item(73, 7).Value = "Zapallo"
$ws.Cells.Item(73, 8).Value = "Camote"
$ws.Cells.Item(73, 9).Value = "1a nueva(o)"
$ws.Cells.Item(73, 10).Value = 300
$ws.Cells.Item(73, 11).Value = 350
$ws.Cells.Item(73, 12).Value = 400
$ws.Cells.Item(73, 13).Value = 375
$ws.Cells.Item(73, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(73, 15).Value = "Región del Maule"
$ws.Cells.Item(73, 16).Value = 375
$ws.Cells.Item(73, 17).Value = 1
$ws.Cells.Item(73, 18).Value = "Hortaliza"
